# Scheduled runner update: refresh scraped market-board price snapshots
# (currentAveragePrice / currentAveragePriceNQ / LevePriceNQ / LeveProfitNQ
# and the derived HQ columns) across the Tonberry_Profits leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1730.5
$ws.Range("I19").Value = 1596.1111
$ws.Range("K19").Value = 1596.1111
$ws.Range("M19").Value = -1421.1111

$ws.Range("H53").Value = 613.3333
$ws.Range("I53").Value = 500
$ws.Range("J53").Value = 636
$ws.Range("K53").Value = 500
$ws.Range("L53").Value = 636
$ws.Range("M53").Value = 137
$ws.Range("N53").Value = -1910

$ws.Range("H58").Value = 2239.2
$ws.Range("I58").Value = 398.66666
$ws.Range("K58").Value = 1195.99998
$ws.Range("M58").Value = -1045.99998

$ws.Range("H62").Value = 8207.375
$ws.Range("I62").Value = 6736
$ws.Range("K62").Value = 6736
$ws.Range("M62").Value = -6112

$ws.Range("H65").Value = 8207.375
$ws.Range("I65").Value = 6736
$ws.Range("K65").Value = 33680
$ws.Range("M65").Value = -30560

$ws.Range("H135").Value = 516.5454999999999
$ws.Range("I135").Value = 409.22223
$ws.Range("K135").Value = 3683.00007
$ws.Range("M135").Value = -1148.00007

$ws.Range("H137").Value = 39221.297
$ws.Range("I137").Value = 984.6667
$ws.Range("J137").Value = 50146.047
$ws.Range("K137").Value = 2954.0001
$ws.Range("L137").Value = 150438.141
$ws.Range("M137").Value = -404.0001000000002
$ws.Range("N137").Value = -155538.141

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2852.0356
$ws.Range("I32").Value = 2014.2046
$ws.Range("J32").Value = 5924.0835
$ws.Range("K32").Value = 2014.2046
$ws.Range("L32").Value = 5924.0835
$ws.Range("M32").Value = -1727.2046
$ws.Range("N32").Value = -6498.0835

$ws.Range("H45").Value = 2737.1936
$ws.Range("I45").Value = 2671.25
$ws.Range("K45").Value = 2671.25
$ws.Range("M45").Value = -2294.25

$ws.Range("H61").Value = 4011.7334
$ws.Range("J61").Value = 4964.75
$ws.Range("L61").Value = 4964.75
$ws.Range("N61").Value = -5388.75

$ws.Range("H74").Value = 1817.6428
$ws.Range("I74").Value = 964.7
$ws.Range("K74").Value = 964.7
$ws.Range("M74").Value = -90.70000000000005

$ws.Range("H77").Value = 1817.6428
$ws.Range("I77").Value = 964.7
$ws.Range("K77").Value = 4823.5
$ws.Range("M77").Value = -455.5

$ws.Range("H132").Value = 3107.647
$ws.Range("I132").Value = 2986.4614
$ws.Range("J132").Value = 3501.5
$ws.Range("K132").Value = 8959.3842
$ws.Range("L132").Value = 10504.5
$ws.Range("M132").Value = -6429.3842
$ws.Range("N132").Value = -15564.5

$ws.Range("H136").Value = 4011.7334
$ws.Range("J136").Value = 4964.75
$ws.Range("L136").Value = 14894.25
$ws.Range("N136").Value = -19994.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 1850
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H134").Value = 3228
$ws.Range("I134").Value = 2911.5715
$ws.Range("J134").Value = 3966.3333
$ws.Range("K134").Value = 8734.7145
$ws.Range("L134").Value = 11898.9999
$ws.Range("M134").Value = -6199.7145
$ws.Range("N134").Value = -16968.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1234.875
$ws.Range("I31").Value = 951.5294
$ws.Range("J31").Value = 1390.258
$ws.Range("K31").Value = 951.5294
$ws.Range("L31").Value = 1390.258
$ws.Range("M31").Value = -656.5294
$ws.Range("N31").Value = -1980.258

$ws.Range("H34").Value = 1234.875
$ws.Range("I34").Value = 951.5294
$ws.Range("J34").Value = 1390.258
$ws.Range("K34").Value = 951.5294
$ws.Range("L34").Value = 1390.258
$ws.Range("M34").Value = -749.5294
$ws.Range("N34").Value = -1794.258

$ws.Range("H107").Value = 483.89474
$ws.Range("I107").Value = 405.52942
$ws.Range("J107").Value = 1150
$ws.Range("K107").Value = 405.52942
$ws.Range("L107").Value = 1150
$ws.Range("M107").Value = 1514.47058
$ws.Range("N107").Value = -4990

$ws.Range("H134").Value = 1972.8125
$ws.Range("I134").Value = 1286.7
$ws.Range("K134").Value = 3860.1
$ws.Range("M134").Value = -1325.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1500088.6
$ws.Range("I4").Value = 1500088.6
$ws.Range("K4").Value = 4500265.800000001
$ws.Range("M4").Value = -4500153.800000001

$ws.Range("H14").Value = 154.25
$ws.Range("I14").Value = 154.25
$ws.Range("K14").Value = 462.75
$ws.Range("M14").Value = -289.75

$ws.Range("H34").Value = 764.8333
$ws.Range("I34").Value = 337.5
$ws.Range("K34").Value = 1012.5
$ws.Range("M34").Value = -928.5

$ws.Range("H61").Value = 199

$ws.Range("H117").Value = 28571712
$ws.Range("I117").Value = 464.66666
$ws.Range("J117").Value = 71428584
$ws.Range("K117").Value = 1393.99998
$ws.Range("L117").Value = 214285752
$ws.Range("M117").Value = 2048.00002
$ws.Range("N117").Value = -214292636

$ws.Range("H131").Value = 20030122
$ws.Range("J131").Value = 35771.145
$ws.Range("L131").Value = 107313.435
$ws.Range("N131").Value = -117393.435

$ws.Range("H137").Value = 4693.9287
$ws.Range("I137").Value = 2460
$ws.Range("J137").Value = 5587.5
$ws.Range("K137").Value = 7380
$ws.Range("L137").Value = 16762.5
$ws.Range("M137").Value = -2280
$ws.Range("N137").Value = -26962.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3455.6667
$ws.Range("I102").Value = 3638.375
$ws.Range("K102").Value = 3638.375
$ws.Range("M102").Value = -2016.375

$ws.Range("H132").Value = 1376644.4
$ws.Range("I132").Value = 2138579.5
$ws.Range("J132").Value = 5161
$ws.Range("K132").Value = 6415738.5
$ws.Range("L132").Value = 15483
$ws.Range("M132").Value = -6413208.5
$ws.Range("N132").Value = -20543

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1354.125
$ws.Range("I46").Value = 866.5
$ws.Range("K46").Value = 866.5
$ws.Range("M46").Value = -678.5

$ws.Range("H132").Value = 2899.394
$ws.Range("I132").Value = 1352.7778
$ws.Range("J132").Value = 3479.375
$ws.Range("K132").Value = 4058.3334
$ws.Range("L132").Value = 10438.125
$ws.Range("M132").Value = -1528.3334
$ws.Range("N132").Value = -15498.125

$ws.Range("H136").Value = 5094.7144
$ws.Range("J136").Value = 6872.1816
$ws.Range("L136").Value = 20616.5448
$ws.Range("N136").Value = -25716.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 83007.30499999999
$ws.Range("I122").Value = 86714.5
$ws.Range("J122").Value = 1449
$ws.Range("K122").Value = 260143.5
$ws.Range("L122").Value = 4347
$ws.Range("M122").Value = -257693.5
$ws.Range("N122").Value = -9247

$ws.Range("H126").Value = 5104.9
$ws.Range("J126").Value = 6255.7144
$ws.Range("L126").Value = 18767.1432
$ws.Range("N126").Value = -23707.1432
